# Append 9 new user_detail rows (rows 22-30) to Sheet1, matching the
# existing table's layout/styling, then update the sheet view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space (U+00A0) used between first/last name in the
# existing data -- keep the same convention for the new rows.
$nbsp = [char]0x00A0

$newRows = @(
    @(110021, 7316931025, "Magdalena${nbsp}Weber",   "magdalena.weber@xyz.com",   932122450),
    @(110022, 9137847236, "Adrienne${nbsp}Hoffman",  "adrienne.hoffman@xyz.com",  848488000),
    @(110023, 8428758532, "Adrienne${nbsp}Mcgee",    "adrienne.mcgee@xyz.com",    894773246),
    @(110024, 9804209494, "Amare${nbsp}Coleman",     "amare.coleman@xyz.com",     956554588),
    @(110025, 7105248214, "Dawson${nbsp}Ibarra",     "dawson.ibarra@xyz.com",     765455583),
    @(110026, 9316557128, "Elvis${nbsp}Mcmillan",    "elvis.mcmillan@xyz.com",    884282274),
    @(110027, 8103486949, "Steve${nbsp}George",      "steve.george@xyz.com",      971073663),
    @(110028, 9601932866, "Colton${nbsp}Elliott",    "colton.elliott@xyz.com",    809908673),
    @(110029, 9317596765, "Carolyn${nbsp}Rodriguez", "carolyn.rodriguez@xyz.com", 818876429)
)

$startRow = 22

# Fill column-by-column (not row-by-row) so new shared-string entries land
# in the same order the source workbook used: all names, then all emails.
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $newRows[$i][0]
}
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $newRows[$i][1]
}
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $newRows[$i][2]
}
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 4).Value = $newRows[$i][3]
}
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = "ACT"
    $ws.Cells.Item($r, 7).Value = "eng"
    $ws.Cells.Item($r, 8).Value = "PWD"
    $ws.Cells.Item($r, 9).Value = $true
    $ws.Cells.Item($r, 10).Value = "superadmin"
    $ws.Cells.Item($r, 11).Value = "now()"
}

# Match the style already used on column I (is_active) in the existing
# rows -- those boolean cells are explicitly left-aligned.
$ws.Range("I22:I30").HorizontalAlignment = -4131

# Update the view: scrolled down + the newly-added range selected.
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("A22:K30").Select()
